$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[1] = @('2025-12-25 16:22', '2025-12-25 14:25', '2025-12-25 14:18', '2025-12-25 14:04', '2025-12-25 12:29', '2025-12-25 10:31', '2025-12-25 09:25', '2025-12-25 04:11', '2025-12-25 02:11', '2025-12-25 00:27', '2025-12-24 22:22', '2025-12-24 20:12', '2025-12-24 19:14', '2025-12-24 16:16', '2025-12-24 14:28', '2025-12-24 12:31', '2025-12-24 10:30', '2025-12-24 09:22', '2025-12-24 04:11', '2025-12-24 02:11', '2025-12-24 00:29', '2025-12-23 22:26', '2025-12-23 20:17', '2025-12-23 19:20', '2025-12-23 16:20', '2025-12-23 14:27', '2025-12-23 12:32', '2025-12-23 10:31', '2025-12-23 09:25', '2025-12-23 04:13', '2025-12-23 02:12', '2025-12-23 00:28', '2025-12-22 22:24', '2025-12-22 20:14', '2025-12-22 19:18', '2025-12-22 16:18', '2025-12-22 14:27', '2025-12-22 12:32', '2025-12-22 10:32', '2025-12-22 09:32', '2025-12-22 04:10', '2025-12-22 02:10', '2025-12-22 00:24', '2025-12-21 22:17', '2025-12-21 20:09', '2025-12-21 19:06', '2025-12-21 16:12', '2025-12-21 14:19', '2025-12-21 12:27', '2025-12-21 10:29', '2025-12-21 09:28', '2025-12-21 04:09', '2025-12-21 02:07', '2025-12-21 00:52', '2025-12-21 00:32', '2025-12-20 23:33', '2025-12-20 23:18', '2025-12-20 22:00', '2025-12-20 13:20')
$rows[2] = @(929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, $null, $null, 929, 929, $null, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929, 929)
$rows[3] = @($null, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, $null, $null, 569, 569, 569, $null, 569, $null, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569)
$rows[4] = @(299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299)
$rows[5] = @(569, $null, 569, 569, 569, 569, 569, 569, $null, $null, 569, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569)
$rows[6] = @($null, 499, 499, 499, 499, 499, 499, 499, 499, $null, 499, $null, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499)
$rows[7] = @(569, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, $null, $null, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569)
$rows[8] = @($null, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, 929, $null, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[9] = @($null, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299)
$rows[10] = @(299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, $null, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, $null, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299)
$rows[11] = @($null, 2997, 2997, 2997, 2997, 2997, 2997, 2997, $null, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, 2997, $null, 2997, 2997, 2997, 2997, 2997, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[12] = @($null, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569)
$rows[13] = @($null, 569, 569, 569, $null, 569, 569, 569, $null, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, 569, $null, 569, $null, 569, 569, 569, 569, 569, 569, $null, 569, 569, 569, $null, 569, 569, 569)
$rows[14] = @(794, 794, 794, 794, 794, 794, 794, 794, 794, $null, 794, 794, $null, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, 794, $null, 794, 794, 794, 499, 499, 499, 499, 499, 499, $null, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499)
$rows[15] = @($null, 499, 499, 499, 499, 499, 499, 499, $null, $null, 499, $null, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, $null, 499, $null, 499, 499, 499, 499, 499, 499, 499)
$rows[16] = @(299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299)
$rows[17] = @($null, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[18] = @(499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, $null, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, $null, $null, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499)
$rows[19] = @(1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, $null, $null, 1299, 1299, $null, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1497, 1497, 1497, $null, $null, 465, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 1497, 2997)
$rows[20] = @($null, 569, 569, 569, 569, 569, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[21] = @(499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, $null, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, $null, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499, 499)
$rows[22] = @($null, 299, 299, 299, 299, 299, 299, 299, $null, $null, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, 299, $null, 299, $null, 299, 299, 299, 299, 299, 299, $null, 299, 299, 299, 299, 299, 299, 299)
$rows[23] = @($null, 1299, 1299, 1299, 1299, 1299, 1299, 1299, $null, 1299, 1299, $null, $null, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, $null, $null, 1299, 1299, 1299, 1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, 1299, 1299, 1299)
$rows[24] = @(569, 569, 569, 569, 569, 569, 929, 929, 929, 929, 929, $null, $null, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, $null, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[25] = @($null, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, $null, $null, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, 929, $null, $null, 929, 929, 929, 929, 929, 929, 929, $null, 929, 929, 929, 929, 929, 929, 929)
$rows[26] = @($null, 1299, 1299, 1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, 1299, $null, 1299, $null, 1299, 1299, 1299, 1299, 1299, 1299, $null, 1299, 1299, 1299, 1299, 1299, 1299, 1299)
$rows[27] = @($null, 199, 199, 199, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rows[28] = @($null, 213, 213, 213, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rows[29] = @($null, 251, 251, 251, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rows[30] = @($null, 231, 231, 231, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rows[31] = @($null, 299, 299, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rows[32] = @(254, 251, 251, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

for ($r = 1; $r -le 32; $r++) {
    $rowvals = $rows[$r]
    for ($i = 0; $i -lt $rowvals.Length; $i++) {
        $col = $i + 2
        $v = $rowvals[$i]
        $cell = $ws.Cells.Item($r, $col)
        if ($v -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

$ws.Columns.Item(60).ColumnWidth = 21